$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-16 (columns E:S) with refreshed computed values ---
# Row 2
$ws.Cells.Item(2, 5).Value = 27779
$ws.Cells.Item(2, 6).Value = 5.733927907840951
$ws.Cells.Item(2, 7).Value = 0.9571558796718381
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 4
$ws.Cells.Item(2, 13).Value = 3
$ws.Cells.Item(2, 14).Value = 3
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = -0.004426163869981436

# Row 3
$ws.Cells.Item(3, 5).Value = 32
$ws.Cells.Item(3, 6).Value = 5.881428122230798
$ws.Cells.Item(3, 7).Value = 2.271346095411742
$ws.Cells.Item(3, 8).Value = 43
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 2.666666666666667
$ws.Cells.Item(3, 13).Value = 2
$ws.Cells.Item(3, 14).Value = 2
$ws.Cells.Item(3, 15).Value = 4.770802919708029
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0.2958637469586374

# Row 4
$ws.Cells.Item(4, 5).Value = 31
$ws.Cells.Item(4, 6).Value = 5.881428122230798
$ws.Cells.Item(4, 7).Value = 2.127013065937418
$ws.Cells.Item(4, 8).Value = 35
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 3.333333333333333
$ws.Cells.Item(4, 13).Value = 1
$ws.Cells.Item(4, 14).Value = 4
$ws.Cells.Item(4, 15).Value = 3.869090909090909
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0.3684848484848485

# Row 5
$ws.Cells.Item(5, 5).Value = 26
$ws.Cells.Item(5, 6).Value = 5.934025097904697
$ws.Cells.Item(5, 7).Value = 2.233363719234288
$ws.Cells.Item(5, 8).Value = 42
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 2.666666666666667
$ws.Cells.Item(5, 13).Value = 2
$ws.Cells.Item(5, 14).Value = 2
$ws.Cells.Item(5, 15).Value = 4.56
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0.2895238095238095

# Row 6
$ws.Cells.Item(6, 5).Value = 28
$ws.Cells.Item(6, 6).Value = 5.934025097904697
$ws.Cells.Item(6, 7).Value = 2.233363719234288
$ws.Cells.Item(6, 8).Value = 42
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 2.666666666666667
$ws.Cells.Item(6, 13).Value = 2
$ws.Cells.Item(6, 14).Value = 2
$ws.Cells.Item(6, 15).Value = 4.592805755395683
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0.2916067146282973

# Row 7
$ws.Cells.Item(7, 5).Value = 29
$ws.Cells.Item(7, 6).Value = 5.994625961181145
$ws.Cells.Item(7, 7).Value = 7.604071710726265
$ws.Cells.Item(7, 8).Value = 39
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 17
$ws.Cells.Item(7, 12).Value = 4
$ws.Cells.Item(7, 13).Value = 2
$ws.Cells.Item(7, 14).Value = 2
$ws.Cells.Item(7, 15).Value = 4.28014440433213
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 0
$ws.Cells.Item(7, 18).Value = 1.865703971119134
$ws.Cells.Item(7, 19).Value = 0.4389891696750903

# Row 8
$ws.Cells.Item(8, 5).Value = 26
$ws.Cells.Item(8, 6).Value = 5.994625961181145
$ws.Cells.Item(8, 7).Value = 7.9003342449104
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 39
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 17
$ws.Cells.Item(8, 12).Value = 4
$ws.Cells.Item(8, 13).Value = 2
$ws.Cells.Item(8, 14).Value = 2
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 4.234285714285714
$ws.Cells.Item(8, 17).Value = 0
$ws.Cells.Item(8, 18).Value = 1.845714285714286
$ws.Cells.Item(8, 19).Value = 0.4342857142857143

# Row 9
$ws.Cells.Item(9, 5).Value = 29
$ws.Cells.Item(9, 6).Value = 5.994625961181145
$ws.Cells.Item(9, 7).Value = 7.9003342449104
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 39
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 17
$ws.Cells.Item(9, 12).Value = 4
$ws.Cells.Item(9, 13).Value = 2
$ws.Cells.Item(9, 14).Value = 2
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 4.28014440433213
$ws.Cells.Item(9, 17).Value = 0
$ws.Cells.Item(9, 18).Value = 1.865703971119134
$ws.Cells.Item(9, 19).Value = 0.4389891696750903

# Row 10
$ws.Cells.Item(10, 5).Value = 29
$ws.Cells.Item(10, 6).Value = 6.247320126918789
$ws.Cells.Item(10, 7).Value = 8.340929808568864
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 42
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 18
$ws.Cells.Item(10, 12).Value = 4
$ws.Cells.Item(10, 13).Value = 2
$ws.Cells.Item(10, 14).Value = 2
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 4.609386281588447
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 1.975451263537906
$ws.Cells.Item(10, 19).Value = 0.4389891696750903

# Row 11
$ws.Cells.Item(11, 5).Value = 14
$ws.Cells.Item(11, 6).Value = 6.334219478032186
$ws.Cells.Item(11, 7).Value = 8.386508659981809
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 43
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 18
$ws.Cells.Item(11, 12).Value = 4
$ws.Cells.Item(11, 13).Value = 2
$ws.Cells.Item(11, 14).Value = 2
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 4.476712328767123
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 1.873972602739726
$ws.Cells.Item(11, 19).Value = 0.4164383561643835

# Row 12
$ws.Cells.Item(12, 5).Value = 31
$ws.Cells.Item(12, 6).Value = 6.35508675642455
$ws.Cells.Item(12, 7).Value = 7.885141294439422
$ws.Cells.Item(12, 8).Value = 1
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 7
$ws.Cells.Item(12, 11).Value = 19
$ws.Cells.Item(12, 12).Value = 6.999999999999997
$ws.Cells.Item(12, 13).Value = 4
$ws.Cells.Item(12, 14).Value = 2
$ws.Cells.Item(12, 15).Value = 0.1105454545454545
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 0.7738181818181817
$ws.Cells.Item(12, 18).Value = 2.100363636363636
$ws.Cells.Item(12, 19).Value = 0.7738181818181814

# Row 13
$ws.Cells.Item(13, 5).Value = 28
$ws.Cells.Item(13, 6).Value = 6.542892261955807
$ws.Cells.Item(13, 7).Value = 7.964904284412076
$ws.Cells.Item(13, 8).Value = 1
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 7
$ws.Cells.Item(13, 11).Value = 19
$ws.Cells.Item(13, 12).Value = 7.33333333333333
$ws.Cells.Item(13, 13).Value = 3
$ws.Cells.Item(13, 14).Value = 3
$ws.Cells.Item(13, 15).Value = 0.1093525179856115
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 0.7654676258992805
$ws.Cells.Item(13, 18).Value = 2.077697841726619
$ws.Cells.Item(13, 19).Value = 0.8019184652278174

# Row 14
$ws.Cells.Item(14, 5).Value = 29
$ws.Cells.Item(14, 6).Value = 6.734413858159677
$ws.Cells.Item(14, 7).Value = 8.508052263749665
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 7
$ws.Cells.Item(14, 11).Value = 20
$ws.Cells.Item(14, 12).Value = 8.33333333333333
$ws.Cells.Item(14, 13).Value = 4
$ws.Cells.Item(14, 14).Value = 3
$ws.Cells.Item(14, 15).Value = 0.1097472924187726
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0.768231046931408
$ws.Cells.Item(14, 18).Value = 2.194945848375451
$ws.Cells.Item(14, 19).Value = 0.9145607701564378

# Row 15
$ws.Cells.Item(15, 5).Value = 28
$ws.Cells.Item(15, 6).Value = 6.812451762284538
$ws.Cells.Item(15, 7).Value = 8.508052263749665
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 7
$ws.Cells.Item(15, 11).Value = 20
$ws.Cells.Item(15, 12).Value = 8.33333333333333
$ws.Cells.Item(15, 13).Value = 3
$ws.Cells.Item(15, 14).Value = 4
$ws.Cells.Item(15, 15).Value = 0.1093525179856115
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 0.7654676258992805
$ws.Cells.Item(15, 18).Value = 2.18705035971223
$ws.Cells.Item(15, 19).Value = 0.9112709832134289

# Row 16
$ws.Cells.Item(16, 5).Value = 28
$ws.Cells.Item(16, 6).Value = 6.812451762284538
$ws.Cells.Item(16, 7).Value = 8.470069887572212
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 7
$ws.Cells.Item(16, 11).Value = 20
$ws.Cells.Item(16, 12).Value = 8.33333333333333
$ws.Cells.Item(16, 13).Value = 4
$ws.Cells.Item(16, 14).Value = 3
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 0.7654676258992805
$ws.Cells.Item(16, 18).Value = 2.18705035971223
$ws.Cells.Item(16, 19).Value = 0.9112709832134289

# --- Add new row 17 for Recluta1 ---
$ws.Cells.Item(17, 1).Value = "Recluta1"
$ws.Cells.Item(17, 2).Value = 44713
$ws.Cells.Item(17, 3).Value = 1.706849315068493
$ws.Cells.Item(17, 4).Value = "Padawan-Sin Fijo"
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 6.812451762284538
$ws.Cells.Item(17, 7).Value = 8.709358857490173
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 7
$ws.Cells.Item(17, 11).Value = 20
$ws.Cells.Item(17, 12).Value = 9.333333333333332
$ws.Cells.Item(17, 13).Value = 4
$ws.Cells.Item(17, 14).Value = 4
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 0.6954248366013072
$ws.Cells.Item(17, 18).Value = 1.986928104575163
$ws.Cells.Item(17, 19).Value = 0.9272331154684095

# Match date style/number format of column B used for other rows (YYYY-MM-DD)
$ws.Cells.Item(17, 2).NumberFormat = $ws.Cells.Item(16, 2).NumberFormat
